$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 42608.893425925926
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 39
$ws.Range("E6").Value = 73
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 16177
$ws.Range("H6").Value = 17712
$ws.Range("I6").Value = 3190
$ws.Range("J6").Value = 354
$ws.Range("K6").Value = 234
$ws.Range("L6").Value = 17
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = "Noun"

$ws.Range("A6").NumberFormat = "m/d/yy h:mm"
